# Fruta / hortaliza, semanal
# Insert 4 new daily-price rows for "Vega Modelo de Temuco" - Palta (Hass)
# right before the existing row 831, shifting the subsequent rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 831 (existing rows 831:904 shift down to 835:908)
$ws.Rows.Item(831).Resize(4).Insert()

# Common / fixed column values shared by every row in this data block
$marketId   = 10
$marketName = "Vega Modelo de Temuco"
$region     = "La Araucanía"
$catId      = 9
$catName    = "Fruta"
$subCatId   = 100106
$subCatName = "Oleaginosos"
$prodId     = 100106002
$prodName   = "Palta"
$variety    = "Hass"

$newRows = @(
    @{ Row = 831; D = 44578; L = "Especial"; M = 180; N = 3500; O = 3800; P = 3708; Q = "`$/kilo (en bandeja de 18 kilos)"; R = "Provincia de Quillota"; S = 3708; T = 1 },
    @{ Row = 832; D = 44578; L = "Primera";  M = 145; N = 3000; O = 3200; P = 3090; Q = "`$/kilo (en bandeja de 18 kilos)"; R = "Provincia de Quillota"; S = 3090; T = 1 },
    @{ Row = 833; D = 44578; L = "Segunda";  M = 235; N = 2000; O = 2500; P = 2234; Q = "`$/kilo (en bandeja de 18 kilos)"; R = "Provincia de Quillota"; S = 2234; T = 1 },
    @{ Row = 834; D = 44578; L = "Tercera";  M = 410; N = 1700; O = 1800; P = 1751; Q = "`$/kilo (en bandeja de 18 kilos)"; R = "Provincia de Quillota"; S = 1751; T = 1 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value  = $marketId
    $ws.Cells.Item($row, 2).Value  = $marketName
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = $catId
    $ws.Cells.Item($row, 6).Value  = $catName
    $ws.Cells.Item($row, 7).Value  = $subCatId
    $ws.Cells.Item($row, 8).Value  = $subCatName
    $ws.Cells.Item($row, 9).Value  = $prodId
    $ws.Cells.Item($row, 10).Value = $prodName
    $ws.Cells.Item($row, 11).Value = $variety
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
